# Fruta / hortaliza, semanal
# Insert a new weekly price-report row above the existing row 31 (shifting
# every subsequent data row down by one, including the previously-last row
# 50 which becomes row 51), then populate the new row 31 with this week's
# reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 31:50 down to 32:51 and open up a blank row 31.
$ws.Rows(31).Insert()

$ws.Range("A31").Value2 = 10
$ws.Range("B31").Value2 = "Vega Modelo de Temuco"
$ws.Range("C31").Value2 = "La Araucanía"
$ws.Range("D31").Value2 = 44452
$ws.Range("E31").Value2 = 9
$ws.Range("F31").Value2 = "Fruta"
$ws.Range("G31").Value2 = 100108
$ws.Range("H31").Value2 = "Tropicales y subtropicales"
$ws.Range("I31").Value2 = 100108004
$ws.Range("J31").Value2 = "Papaya"
$ws.Range("K31").Value2 = "Cultivar IV Región"
$ws.Range("L31").Value2 = "Primera"
$ws.Range("M31").Value2 = 65
$ws.Range("N31").Value2 = 21000
$ws.Range("O31").Value2 = 21000
$ws.Range("P31").Value2 = 21000
$ws.Range("Q31").Value2 = "$/bandeja 10 kilos"
$ws.Range("R31").Value2 = "Provincia del Elquí"
$ws.Range("S31").Value2 = 2100
$ws.Range("T31").Value2 = 10
